$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "39.258.48"
$ws.Range("E2").Value = "  -1.54%  "

Set-TextValue "D3" "2.189.58"
$ws.Range("E3").Value = "  -6.08%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D5" "295.66"
$ws.Range("E5").Value = "  -3.77%  "

Set-TextValue "D6" "81.35"
$ws.Range("E6").Value = "  -3.04%  "

$ws.Range("E7").Value = "  -3.75%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -3.62%  "

Set-TextValue "D10" "0.0768"
$ws.Range("E10").Value = "  -5.95%  "

Set-TextValue "D11" "29.05"
$ws.Range("E11").Value = "  -2.72%  "

Set-TextValue "D12" "47.09"
$ws.Range("E12").Value = "  -10.24%  "

$ws.Range("E13").Value = "  -2.15%  "

Set-TextValue "D14" "6.23"
$ws.Range("E14").Value = "  -2.37%  "

Set-TextValue "D15" "2.526.51"
$ws.Range("E15").Value = "  -6.56%  "

Set-TextValue "D16" "13.89"
$ws.Range("E16").Value = "  -5.92%  "

Set-TextValue "D17" "2.192.67"
$ws.Range("E17").Value = "  -7.14%  "

$ws.Range("E18").Value = "  -5.12%  "

Set-TextValue "D19" "39.126.08"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("E20").Value = "  -3.43%  "

Set-TextValue "D21" "5.68"
$ws.Range("E21").Value = "  -6.00%  "

Set-TextValue "D22" "64.74"
$ws.Range("E22").Value = "  -4.35%  "

Set-TextValue "D23" "10.17"
$ws.Range("E23").Value = "  -3.59%  "

Set-TextValue "D24" "224.85"
$ws.Range("E24").Value = "  -3.76%  "

$ws.Range("E25").Value = "  +0.03%  "

Set-TextValue "D26" "2.39"
$ws.Range("E26").Value = "  -6.00%  "

Set-TextValue "D27" "1.79"
$ws.Range("E27").Value = "  +0.42%  "

Set-TextValue "D28" "22.43"
$ws.Range("E28").Value = "  -3.62%  "

$ws.Range("E29").Value = "  -1.63%  "

Set-TextValue "D30" "9.03"
$ws.Range("E30").Value = "  -1.31%  "

Set-TextValue "D31" "150.10"
$ws.Range("E31").Value = "  -0.55%  "

Set-TextValue "D32" "31.40"
$ws.Range("E32").Value = "  -7.92%  "

$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("E34").Value = "  -6.27%  "

$ws.Range("E35").Value = "  -3.69%  "

Set-TextValue "D36" "0.0688"
$ws.Range("E36").Value = "  -4.62%  "

$ws.Range("E37").Value = "  -3.32%  "

# Row 38 and 39 swap coin identities (Celestia <-> Kaspa) with new price/volume data
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.0960"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D39" "15.22"
$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("E40").Value = "  -5.09%  "

$ws.Range("E41").Value = "  -3.48%  "

Set-TextValue "D42" "3.58"
$ws.Range("E42").Value = "  -5.28%  "

Set-TextValue "D43" "1.887.76"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("E44").Value = "  -11.48%  "

$ws.Range("E45").Value = "  -2.27%  "

Set-TextValue "D46" "16.09"
$ws.Range("E46").Value = "  -7.94%  "

$ws.Range("E47").Value = "  -4.53%  "

Set-TextValue "D48" "2.59"
$ws.Range("E48").Value = "  -1.54%  "

Set-TextValue "D49" "2.398.27"
$ws.Range("E49").Value = "  -7.10%  "

Set-TextValue "D50" "70.68"
$ws.Range("E50").Value = "  +0.09%  "

Set-TextValue "D51" "86.59"
$ws.Range("E51").Value = "  -5.77%  "
